# C5-PowerPoint.pptx — commit "Tue, May 05, 2020 12:07:12 PM"
#
# Change: the table on Slide 6 ("SOURCES OF FINANCE") had its table style
# (Table Design gallery pick) changed from style {380883B5-BCB4-43F9-8120-
# 36C71230EE34} to style {22CF40C7-503A-4801-AAAC-F93D3F6BC8E4}.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(6)
$sh = $s.Shapes.Item(2)
$tbl = $sh.Table

$tbl.ApplyStyle("{22CF40C7-503A-4801-AAAC-F93D3F6BC8E4}")
